{"js": "// Replace the date paragraph's text (\"2024-05-30\") with \"Invalid Date\",\n// split across three runs (\"Invalid\", \" \", \"Date\") to mirror the\n// word-by-word run layout used elsewhere in the document.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/style,items/text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].style === \"Date\") {\n    target = paragraphs.items[i];\n    break;\n  }\n}\nif (!target) {\n  throw new Error(\"Could not find a paragraph with style 'Date'.\");\n}\n\n// Flat-OPC wrapped OOXML so the paragraph's runs come in as three\n// discrete <w:r> elements instead of being coalesced into one.\nconst ooxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:pPr><w:pStyle w:val=\"Date\"/></w:pPr>' +\n  '<w:r><w:t xml:space=\"preserve\">Invalid</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\">Date</w:t></w:r>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\ntarget.insertOoxml(ooxml, \"Replace\");\nawait context.sync();\n", "ps1": "# Replace the date paragraph's text (\"2024-05-30\") with \"Invalid Date\",\n# split across three runs (\"Invalid\", \" \", \"Date\") to mirror the\n# word-by-word run layout used elsewhere in the document.\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Style.NameLocal -eq \"Date\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not find a paragraph with style 'Date'.\"\n}\n\n$xml = \"<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>\" +\n       \"<w:pPr><w:pStyle w:val='Date'/></w:pPr>\" +\n       \"<w:r><w:t xml:space='preserve'>Invalid</w:t></w:r>\" +\n       \"<w:r><w:t xml:space='preserve'> </w:t></w:r>\" +\n       \"<w:r><w:t xml:space='preserve'>Date</w:t></w:r>\" +\n       \"</w:p>\"\n\n$target.Range.InsertXML($xml)\n"}
